# Update the "想去人数" (number of people interested) values on the
# "展览" sheet and the aggregated "全部类型" sheet to reflect the latest
# scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 188
$wsExhibit.Range("F3").Value = 505
$wsExhibit.Range("F4").Value = 27
$wsExhibit.Range("F9").Value = 91
$wsExhibit.Range("F10").Value = 1939
$wsExhibit.Range("F11").Value = 8

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 188
$wsAll.Range("F4").Value = 505
$wsAll.Range("F5").Value = 27
$wsAll.Range("F10").Value = 91
$wsAll.Range("F11").Value = 1939
$wsAll.Range("F12").Value = 8
